$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Marking" row's Right-answer total (B11) and the
# "Total" row's Right-answer total (B12) and corrected/total marks (E12)
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 80
$ws.Range("E12").Value = "80/140"
